$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (existing row, target cluster changes from FAPs to ECs)
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.5587383333333333
$ws.Range("H2").Value = 1.676215
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.178057
$ws.Range("N2").Value = 0.534171
$ws.Range("O2").Value = 0.01072706596418309
$ws.Range("P2").Value = 0.01486309420510698
$ws.Range("Q2").Value = 0.09948727141833333
$ws.Range("R2").Value = 0.8953854427649999
$ws.Range("S2").Value = 0.01072706596418309
$ws.Range("T2").Value = 0.01486309420510698

# Update row 3 (existing row, target cluster changes from sCs to FAPs)
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.5587383333333333
$ws.Range("H3").Value = 1.676215
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.509764
$ws.Range("N3").Value = 7.529292
$ws.Range("O3").Value = 0.1512010422647355
$ws.Range("P3").Value = 0.2094995353430986
$ws.Range("Q3").Value = 1.40230135442
$ws.Range("R3").Value = 12.62071218978
$ws.Range("S3").Value = 0.1512010422647355
$ws.Range("T3").Value = 0.2094995353430986

# New row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Lama1"
$ws.Range("C4").Value = "Itgb8"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5587383333333333
$ws.Range("H4").Value = 1.676215
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05389233333333333
$ws.Range("N4").Value = 0.161677
$ws.Range("O4").Value = 0.003246750280137314
$ws.Range("P4").Value = 0.004498597793214309
$ws.Range("Q4").Value = 0.03011171250611111
$ws.Range("R4").Value = 0.271005412555
$ws.Range("S4").Value = 0.003246750280137314
$ws.Range("T4").Value = 0.004498597793214309

# New row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Lama1"
$ws.Range("C5").Value = "Itgb8"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5587383333333333
$ws.Range("H5").Value = 1.676215
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 13.8571405
$ws.Range("N5").Value = 27.714281
$ws.Range("O5").Value = 0.8348251414909441
$ws.Range("P5").Value = 0.7711387726585802
$ws.Range("Q5").Value = 7.742515587735833
$ws.Range("R5").Value = 46.455093526415
$ws.Range("S5").Value = 0.8348251414909441
$ws.Range("T5").Value = 0.7711387726585802
